# Data Collection.xlsx -- "added 6 additional products with links" edit
#
# Adds 6 new products (Hogwarts Legacy, Jedi Survivor, MarioKart8 Deluxe,
# Ninja Woodfire Pizza Oven, Igloo 24qt Cooler, Dyson Ball Animal 3 vacuum)
# to the productCollectionList sheet, and their target/amazon/bestbuy URLs
# to the productDataCapture sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("productCollectionList")
$ws2 = $wb.Worksheets.Item("productDataCapture")

# ---------------------------------------------------------------------
# 1. productCollectionList: new product names in B16:B21
# ---------------------------------------------------------------------
$ws1.Range("B16").Value = "Hogwarts Legacy (XBOX)"
$ws1.Range("B17").Value = "Jedi Survivor PS5"
$ws1.Range("B18").Value = "MarioKart8 Deluxe - Switch"
$ws1.Range("B19").Value = "Ninja Woodfire Pizza Oven"
$ws1.Range("B20").Value = "Igloo 24qt Cooler"
$ws1.Range("B21").Value = "Dyson Ball Animal 3 Extra Upright Vacuum Cleaner"

# ---------------------------------------------------------------------
# 2. productDataCapture: target.com links, rows 16-21 (column D)
# ---------------------------------------------------------------------
$ws2.Range("D16").Value = "https://www.target.com/p/hogwarts-legacy-xbox-series-x/-/A-86267032"
$ws2.Range("D17").Value = "https://www.target.com/p/star-wars-jedi-survivor-playstation-5/-/A-88183520"
$ws2.Range("D18").Value = "https://www.target.com/p/mario-kart-8-deluxe-nintendo-switch/-/A-52161278"
$ws2.Range("D19").Value = "https://www.target.com/p/ninja-woodfire-pizza-oven-8-in-1-outdoor-oven-5-pizza-settings-smoker-ninja-woodfire-technology-electric-oo101/-/A-89401843"
$ws2.Range("D20").Value = "https://www.target.com/p/igloo-imx-24qt-cooler-meteorite/-/A-87378961"
$ws2.Range("D21").Value = "https://www.target.com/p/dyson-ball-animal-3-extra-upright-vacuum/-/A-85269264"

# ---------------------------------------------------------------------
# 3. productDataCapture: amazon.com links, rows 44-49 (column D)
# ---------------------------------------------------------------------
$ws2.Range("D44").Value = "https://www.amazon.com/Hogwarts-Legacy-Xbox-English-Version-Region/dp/B0B9YSJC5S/ref=sr_1_2?crid=1K69XYMHWNR56&dib=eyJ2IjoiMSJ9.Kg0XV9Kot63z1u2KagkEvpCVUtkbGXvapqAJX7afkiu3BdvCRxuMMxHq3-R-Vq88oqItB-jjTuSzsA4OUT3gHqxnGxm_W3ax9Bt1-kNDhIIrfk2DMuaREs1iW3ZpC41UpXGGXpnKaQWx_QR6PweYu_E0fgQaz6JOeKQrHgmpoGnZNkaPqUf-EzAWuLvi4LuOAF_AT5tsuZUIeboaBcrdGFpzu4TWaNBwkw4WmSvH7Vo.uFmuTgI86APUwJ8qSsprJqvZ4Dnr_EwLXhTo2K61-P0&dib_tag=se&keywords=hogwarts+legacy+xbox&qid=1710278855&sprefix=hogwarts+legacy+xbox%2Caps%2C129&sr=8-2"
$ws2.Range("D45").Value = "https://www.amazon.com/Star-Wars-Jedi-Survivor-Playstation-5/dp/B0BPD4JMWT/ref=sr_1_1?crid=1YEWINBQP33M2&dib=eyJ2IjoiMSJ9.xmgn1bEBQh8fF6vN-_ijDTx3NnKfAvuz4KzU8chGgsL67cKZ-GVbOSb0nGWeW8mxQpNntErtqMVDRP2_Qvv6NVonWyAdXrQOG1vF3pQuiK5vfmslLzeTXkEtyMSwY0hE5dvkCY4Udzheb9VLdQgjV4RNjelCbXWMWS3HDuqkxw8RhbOBcyBB1KuJpYwpU0LO3JbTFOrTQoCecarfSPAyc01kEOKHqvHzXAHZPHjTRpo.vqAgJHJ_EuyQ0aHKHSEiLA9mvUUFdQtDGBvQua34gew&dib_tag=se&keywords=jedi+survivor+ps5&qid=1710278950&sprefix=jedi+survivor+ps%2Caps%2C174&sr=8-1"
$ws2.Range("D46").Value = "https://www.amazon.com/Mario-Kart-8-Deluxe-Nintendo-Switch/dp/B01N1037CV/ref=sr_1_1?crid=1PYFQU1N4B68F&dib=eyJ2IjoiMSJ9.8TZU_Of0o7Vk3YWNjDE9oY9I_W3GMrGxHT454fa05jnyQhIJ1QoKGzaub95DXly85ZIxHf3U7vByBD9UHUHObxJwp-_uYmPTubto19uqmmQglLhYd3W-DOvMZz4lfv1l9ezRsNuuryHtJljvddmeIZdaJWrxpna-VU5FHp3braflTwk_ypjvU7neyKr_vIay2KHtiMOKVjqNxAnzmGjWABX-bq-mvp-ghJ0Hq7eY9h4.vohI6h3tcX_pnwUMLahacLD3yr74bAsNKv9eeH5rAp0&dib_tag=se&keywords=mario+kart+8+deluxe&qid=1710279189&sprefix=mario+kart+%2Caps%2C167&sr=8-1"
$ws2.Range("D47").Value = "https://www.amazon.com/Ninja-Woodfire-Technology-terracotta-OO101/dp/B0C6BQVDX3/ref=sr_1_3?crid=3DXGBNWLVTJHJ&dib=eyJ2IjoiMSJ9.0WxsIb1exp9JW6Y3VWlsjnd4AbzdLgvTZ3dzGlCYh7z1OEIvxums8Mey72f26ju-BWrkrofq9sMNvCk49njaQbu-SN8ItGFSBFcq-t91Hsb6AE3rPWtoewiMwvt1h2JtkHhsZTKkYeLo0Y9nnD_eJ4STJOza34YIVUe9YB7S8nTqDnhXuTlMQ4TE6pHn3CnUHoVNJGIlgiU2ncUCADPEKpW5C2f1-rC8jpiw7k3qXsE.uqci_nwiB33tgkfGHQRjqPWrYYjoEoVkB-fSO930oYQ&dib_tag=se&keywords=Ninja%2BWoodfire%2BPizza%2BOven%2C%2B8-in-1%2BOutdoor%2BOven&qid=1710280162&sprefix=ninja%2Bwoodfire%2Bpizza%2Boven%2C%2B8-in-1%2Boutdoor%2Boven%2Caps%2C183&sr=8-3&th=1"
$ws2.Range("D48").Value = "https://www.amazon.com/dp/B08YS7QZLL/ref=twister_B09V6MDSN6?_encoding=UTF8&psc=1"
$ws2.Range("D49").Value = "https://www.amazon.com/Dyson-Animal-Extra-Upright-Vacuum/dp/B0B57QV5TT/ref=sr_1_3?crid=2IQWXZR08I144&dib=eyJ2IjoiMSJ9.t3iEOoxEJUwLPVoSUcRptqioH5l2EvX8s4APBjqpAjI4e18xB7pmNY8TyLEanhvfvEQqqscP3Xi13oKWZQuQ1MaGTAAh82rbXWBT5nZjBYQqPFp7PUAvYcMk9RslrAyAZS0Nwv3r6RQ8M2iBSK7DQAZmAK0ycSydB1nB1cVDRVBJ0bgIwOtJPiE6ohrQceFbK44RLbkqiNcgSv_YM7QC30n8grFBmSgixz1fAwZtEd8.ZIHRSLxkS_22b7rx1F4ymcSJXsCEuBy_aVi0A1DzCpU&dib_tag=se&keywords=Dyson%2B-%2BBall%2BAnimal%2B3%2BExtra%2BUpright&qid=1710280601&sprefix=dyson%2B-%2Bball%2Banimal%2B3%2Bextra%2Bupright%2Caps%2C176&sr=8-3&th=1"

# ---------------------------------------------------------------------
# 4. productDataCapture: bestbuy.com links, new rows 75-80
#    (C column = "bestbuy" site marker, D column = url)
#    NOTE: shared-string insertion order follows the original authoring
#    order (Hogwarts, Jedi, Mario, Igloo, Dyson, Ninja last) even though
#    the rows themselves are laid out Hogwarts..Dyson with Ninja in the
#    4th row -- so D79/D80 are written before D78 to reproduce that.
# ---------------------------------------------------------------------
$ws2.Range("C75").Value = "bestbuy"
$ws2.Range("D75").Value = "https://www.bestbuy.com/site/hogwarts-legacy-standard-edition-xbox-series-x/6502627.p?skuId=6502627"

$ws2.Range("C76").Value = "bestbuy"
$ws2.Range("D76").Value = "https://www.bestbuy.com/site/star-wars-jedi-survivor-standard-edition-playstation-5/6528657.p?skuId=6528657"

$ws2.Range("C77").Value = "bestbuy"
$ws2.Range("D77").Value = "https://www.bestbuy.com/site/mario-kart-8-deluxe-nintendo-switch/5723304.p?skuId=5723304"

$ws2.Range("C79").Value = "bestbuy"
$ws2.Range("D79").Value = "https://www.bestbuy.com/site/igloo-imx-24-quart-cooler-gray/6523836.p?skuId=6523836"

$ws2.Range("C80").Value = "bestbuy"
$ws2.Range("D80").Value = "https://www.bestbuy.com/site/dyson-ball-animal-3-extra-upright-vacuum-with-5-accessories-copper-silver/6504485.p?skuId=6504485"

$ws2.Range("C78").Value = "bestbuy"
$ws2.Range("D78").Value = "https://www.bestbuy.com/site/ninja-woodfire-pizza-oven-8-in-1-outdoor-oven-5-pizza-settings-700f-smoker-woodfire-technology-electric-terracotta-red/6548901.p?skuId=6548901"
$ws2.Hyperlinks.Add($ws2.Range("D78"), $ws2.Range("D78").Value) | Out-Null

# ---------------------------------------------------------------------
# 5. Best-effort view/window state to mirror the re-saved workbook
# ---------------------------------------------------------------------
$ws1.Select()
$excel.ActiveWindow.Zoom = 115
$ws1.Range("C23").Select()

$ws2.Select()
$excel.ActiveWindow.FreezePanes = $false
$ws2.Range("A12").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("D24").Select()
